# "Fix imports and record format"
#
# The data-import sheets were named "Logs"/"Users" (capitalized) but the
# underlying CSV-driven record format/tooling expects lower-case sheet
# names ("logs"/"users"). Rename both sheets; Excel automatically keeps
# the workbook-level defined names ("sheet1"/"sheet2") pointing at the
# renamed sheets.
$wb = $excel.ActiveWorkbook

$logsSheet  = $wb.Worksheets.Item("Logs")
$usersSheet = $wb.Worksheets.Item("Users")

$logsSheet.Name  = "logs"
$usersSheet.Name = "users"

# The previously-active tab was "Users" (now "users"); the fixed import
# makes "logs" (the first tab) the active/selected sheet instead.
$logsSheet.Activate()
